$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $rng = $d.Paragraphs($paraIndex).Range
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2)
    return $ok
}

# 1. Title line: "QA Test Analyst | Chennai, India " -> "Senior Quality Engineer | Chennai, India "
Replace-InParagraph 2 "QA Test Analyst" "Senior Quality Engineer" | Out-Null

# 2. Summary bullet: years of experience wording
Replace-InParagraph 8 "3 years & 7-month experience in the Warehouse Management System." "4 years' experience in the Warehouse Management System." | Out-Null

# 3. Summary bullet: Automated testing tools
Replace-InParagraph 9 "Have Knowledge in Automated testing in Playwright, Selenium, Appium" "Experienced in Automated testing using Playwright with JavaScript and TypeScript" | Out-Null

# 4. Summary bullet: Redis Cache / Salesforce knowledge wording
Replace-InParagraph 10 "Have Knowledge in Redis Cache, Salesforce" "Experienced in Redis Cache, Salesforce" | Out-Null

# 5. Technical summary: Programming Languages
Replace-InParagraph 23 "Java Script, Java" "JavaScript, TypeScript" | Out-Null

# 6. Technical summary: Automation Tool
Replace-InParagraph 24 "Playwright, Selenium, Appium" "Playwright" | Out-Null

# 7. Technical summary: Framework (drop TestNG)
Replace-InParagraph 25 "bject Model design, TestNG" "bject Model design" | Out-Null

# 8. Role lines across the 4 project blocks: append ", Automation Tester"
Replace-InParagraph 37 "Functional Tester" "Functional Tester, Automation Tester" | Out-Null
Replace-InParagraph 44 "Functional Tester" "Functional Tester, Automation Tester" | Out-Null
Replace-InParagraph 51 "Functional Tester" "Functional Tester, Automation Tester" | Out-Null
Replace-InParagraph 58 "Functional Tester" "Functional Tester, Automation Tester" | Out-Null

# 9. Roles and Responsibilities: Test plans / cases / scripts
Replace-InParagraph 67 "Designed and developed Test plans and Test cases." "Designed and developed Test plans, Test cases, Test scripts." | Out-Null

# 10. Roles and Responsibilities: Closely and independently ... (merge runs, same text)
Replace-InParagraph 75 "Worked Closely and independently with various teams, including development and management." "Worked Closely and independently with various teams, including development and management." | Out-Null

# 11. Education table: Bannari Amman Institute (remove spell-check split, merge runs)
Replace-InParagraph 85 "Bannari Amman Institute of Technology, Erode" "Bannari Amman Institute of Technology, Erode" | Out-Null
